# Update "想去人数" (attendance count) values in the "展览" (sheet1) and
# "全部类型" (sheet4) worksheets. Both sheets list the same events, just with
# different row offsets, so the edits are applied twice with a row shift.

$wb = $excel.ActiveWorkbook

# Sheet "展览": row number -> new value for column F
$exhibitionUpdates = @{
    6  = 9526
    7  = 856
    10 = 1687
    11 = 155
    12 = 104
    15 = 446
    18 = 1310
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# Sheet "全部类型": same events, rows shifted by +1
$allTypesUpdates = @{
    7  = 9526
    8  = 856
    11 = 1687
    12 = 155
    13 = 104
    16 = 446
    19 = 1310
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
